$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 new rows before the old "Total" row (row 15), pushing
#    the totals/summary block down to rows 19-25.
# ------------------------------------------------------------------
$ws.Rows("15:18").Insert() | Out-Null

# ------------------------------------------------------------------
# 2. Clone formatting for the 4 new data rows from the existing
#    member-row pattern (row 12), for every column A:X.
# ------------------------------------------------------------------
$ws.Range("A12:X12").Copy() | Out-Null
$ws.Range("A15:X18").PasteSpecial(-4122) | Out-Null

# The bottom-most member row uses the special "corner" style (27) in
# column A. Grab it from the current row 14 (still carrying it) before
# row 14 is reset to the regular style, then stamp it onto row 18 -
# the new bottom-most member row.
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

# Row 14 is no longer the last member row, so it goes back to the
# regular style (matching rows 12/13).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Populate the 4 new member rows (15-18).
# ------------------------------------------------------------------
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "sdivhdsiof"
$ws.Range("C15").Value = "sdifndsif"
$ws.Range("D15").Value = 0
$ws.Range("E15").Formula = "=E11"
$ws.Range("G15").Value = "io"
$ws.Range("H15").Value = 1
$ws.Range("T15").Value = 1

$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "djahit"
$ws.Range("C16").Value = "ihd"
$ws.Range("D16").Value = 0
$ws.Range("E16").Formula = "=E11"
$ws.Range("F16").Value = 1
$ws.Range("T16").Value = 1

$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "okkk"
$ws.Range("C17").Value = "ihd"
$ws.Range("D17").Value = 0
$ws.Range("E17").Formula = "=E11"
$ws.Range("F17").Value = 1
$ws.Range("T17").Value = 1

$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "okkk"
$ws.Range("C18").Value = "ihd"
$ws.Range("D18").Value = 0
$ws.Range("E18").Formula = "=E11"
$ws.Range("F18").Value = 1
$ws.Range("T18").Value = 1

# ------------------------------------------------------------------
# 4. Re-point the "Total" row (now row 19) SUM ranges so they cover
#    the full member block D12:D18 (etc.) instead of D12:D14.
# ------------------------------------------------------------------
$ws.Range("D19").Formula = "=SUM(D12:D18)"
$ws.Range("E19").Formula = "=SUM(E12:E18)"
$ws.Range("F19").Formula = "=SUM(F12:F18)"
$ws.Range("H19").Formula = "=SUM(H12:H18)"
$ws.Range("I19").Formula = "=SUM(I12:I18)"
$ws.Range("J19").Formula = "=SUM(J12:J18)"
$ws.Range("K19").Formula = "=SUM(K12:K18)"
$ws.Range("L19").Formula = "=SUM(L12:L18)"
$ws.Range("M19").Formula = "=SUM(M12:M18)"
$ws.Range("N19").Formula = "=SUM(N12:N18)"
$ws.Range("O19").Formula = "=SUM(O12:O18)"
$ws.Range("P19").Formula = "=SUM(P12:P18)"
$ws.Range("Q19").Formula = "=SUM(Q12:Q18)"
$ws.Range("R19").Formula = "=SUM(R12:R18)"
$ws.Range("S19").Formula = "=SUM(S12:S18)"
$ws.Range("T19").Formula = "=SUM(T12:T18)"
$ws.Range("U19").Formula = "=SUM(U12:U18)"
$ws.Range("V19").Formula = "=SUM(V12:V18)"
$ws.Range("W19").Formula = "=SUM(W12:W18)"

# ------------------------------------------------------------------
# 5. Selection ends up on the new total row's W cell.
# ------------------------------------------------------------------
$ws.Range("W19").Select() | Out-Null
